$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("경찰 불러주세요.", "Call the police."),
    @("불이야!", "Fire!"),
    @("도둑이야!", "Thief!"),
    @("도와주세요!", "Help!")
)

$startRow = 98
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}
